$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 880 (old rows 880..978 shift down to 882..980)
$ws.Rows.Item(880).Insert()
$ws.Rows.Item(880).Insert()

# New row 880: duplicate of former row 880 (Pintón) with an updated sampling date and prices
$row880 = New-Object 'object[,]' 1,20
$row880[0,0]  = 4
$row880[0,1]  = "Feria Lagunitas de Puerto Montt"
$row880[0,2]  = "Los Lagos"
$row880[0,3]  = 45212
$row880[0,4]  = 10
$row880[0,5]  = "Fruta"
$row880[0,6]  = 100108
$row880[0,7]  = "Tropicales y subtropicales"
$row880[0,8]  = 100108006
$row880[0,9]  = "Plátano"
$row880[0,10] = "Sin especificar"
$row880[0,11] = "Pintón"
$row880[0,12] = 600
$row880[0,13] = 25000
$row880[0,14] = 25000
$row880[0,15] = 25000
$row880[0,16] = "`$/caja 20 kilos"
$row880[0,17] = "Ecuador"
$row880[0,18] = 1250
$row880[0,19] = 20
$ws.Range("A880:T880").Value2 = $row880

# New row 881: duplicate of former row 881 (Primera Pintón) with an updated sampling date and prices
$row881 = New-Object 'object[,]' 1,20
$row881[0,0]  = 4
$row881[0,1]  = "Feria Lagunitas de Puerto Montt"
$row881[0,2]  = "Los Lagos"
$row881[0,3]  = 45212
$row881[0,4]  = 10
$row881[0,5]  = "Fruta"
$row881[0,6]  = 100108
$row881[0,7]  = "Tropicales y subtropicales"
$row881[0,8]  = 100108006
$row881[0,9]  = "Plátano"
$row881[0,10] = "Sin especificar"
$row881[0,11] = "Primera Pintón"
$row881[0,12] = 1200
$row881[0,13] = 26000
$row881[0,14] = 27000
$row881[0,15] = 26500
$row881[0,16] = "`$/caja 20 kilos"
$row881[0,17] = "Ecuador"
$row881[0,18] = 1325
$row881[0,19] = 20
$ws.Range("A881:T881").Value2 = $row881
